$d = $word.ActiveDocument
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mo="http://schemas.microsoft.com/office/mac/office/2008/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:mv="urn:schemas-microsoft-com:mac:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 wp14"><w:body>
    <w:p>
      <w:r>
        <w:t>Story</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>This is the end…</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>A meteor strike has destined this planet for doom. But there is some hope. From the ashes, humanity has a chance of rebuilding. The government had been building underground bunkers for decades in preparation for a catastrophic global event like a nuclear fallout, but it turns out the biggest threat to humanity came form above. You have been tasked of defending our city for as long as possible. Buy the citizens enough time to load into the bunkers. Don’t let the meteors hit the streets.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>I need y</w:t>
      </w:r>
      <w:r>
        <w:t>ou to know soldier. You and you platoon will not make it through this</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. This will be your last contribution and the most important. May your lives not be lost in </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>vain.</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>How to play</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Shoot down the meteors using the turret y</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">ou have been assigned with the </w:t>
      </w:r>
      <w:r>
        <w:t>right</w:t>
      </w:r>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> l</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">eft arrows </w:t>
      </w:r>
      <w:r>
        <w:t>and SPACE keys. As you shoot down more meteors you will collect fusion cores that fall from the meteors once they are destroyed. Use these fusion cores to power surrounding turrets by pressing U. Once you collect enough fusion cores yo</w:t>
      </w:r>
      <w:r>
        <w:t>u will</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> be able to upgrade turrets also by pressing U. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Different types of meteors have different strengths and weaknesses. </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Game over</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Score:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Thanks to your valiant effort people have been saved. While billions will still die across the globe, you and your platoon have given us hope. May you rest in peace </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>soldier.</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve">  You will not be forgotten in the new world. </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>How to play (detailed)</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Upgrades</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Once</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> you receive enough fusion cores a randomly chosen turret base’s upgrade button will be lit up. Press U to activate the upgrade. You have no control on which turret to perform the upgrade on. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> The initially upgraded </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>guns</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> includi</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">ng the one the user starts with, </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">does 1 damage per bullet. However the second upgrade </w:t>
      </w:r>
      <w:r>
        <w:t>does 2 damage.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Guns</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>As the player you only have full control of one gun. You can control that gun using the RIGHT and LEFT arrows and press the SPACE bar to fire. There is a rate of fire control, meaning you can’t blanket the entire screen with your bullets.  Your bullets become faster as you fire at meteors closer to the ground. The other guns (</w:t>
      </w:r>
      <w:r>
        <w:t>once</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> you have upgraded them) will fire by themselves however they like. </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>they</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> aren’t the most accurate. You can also upgrade your own gun.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Currency</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> and scoring </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">The user controlled gun gains more fusion and points for shooting down a meteor than an automatic gun. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">30 fusion cores for initial upgrade. 60 fusion cores for further upgrade (including user controlled gun). At the end of the game your left over fusion cores are </w:t>
      </w:r>
      <w:r>
        <w:t>multiplied by 10</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> and </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">added to your score. </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Meteors</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Meteor 1 has 2</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:t xml:space="preserve"> health</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> and goes at an average speed compared to the rest</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. Meteor 2 (orange) has 3 </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">health and travels slow when compared to the rest. Meteor 3 </w:t>
      </w:r>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:t>yellow) has 1 health</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> and travels fast when compared to the rest. Meteor 4 has 8 health and travels very slowly when compared to the rest. </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Objective and lives</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>The objective of the game is to prevent a meteor from hitting the ground. If that happens then the game is over. Each building has a different amount of meteor hits it can take. The tall and medium buildings can take 3 hits before the 4</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="superscript"/>
        </w:rPr>
        <w:t>th</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> meteor hits the ground. The small buildings can only take 2 hits before the 3</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="superscript"/>
        </w:rPr>
        <w:t>rd</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> meteor hits the ground. </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    </w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$c = $d.Content
$c.Text = ""
$c.InsertXML($xml)
Write-Host "Paragraphs count:" $d.Paragraphs.Count
